$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = 41

$ws.Range("F12").Value = 250
$ws.Range("H12").Value = 250

$ws.Range("F14").Value = 61
$ws.Range("H14").Value = 61

$ws.Range("E17").Value = 89

$ws.Range("E27").Value = 306
$ws.Range("F27").Value = 149
$ws.Range("H27").Value = 149

$ws.Range("E47").Value = 420
